# ----------------------------------------------------------------------------
# Add demand-reaction data to table S1:
#   - new worksheet 'AGORA model demand reactions' appended after the existing
#     two sheets (becomes the active tab, matching the source workbook)
#   - four section-header rows (merged A:D, left/center aligned) each followed
#     by its list of demand/sink reactions (id in col A, description in col B)
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# Insert the new worksheet after the last existing sheet, matching tab order
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "AGORA model demand reactions"

# Section-header rows: label in column A, merged & centred across A:D
$headers = @(
    @(1, 'AGORA str. NAP07 GENRE demand reactions:'),
    @(12, 'AGORA str. NAP08 demand reactions:'),
    @(23, 'AGORA str. CD196 GENRE demand reactions:'),
    @(45, 'AGORA str. R20291 GENRE demand reactions:')
)

# Data rows: reaction id in column A, description in column B
$data = @(
    @(2, 'DM_4HBA', 'Sink needed to allow 4-hydroxy-benzoate to leave system'),
    @(3, 'DM_HQN', 'Sink needed to allow Hydroquinone to leave system'),
    @(4, 'DM_btn', 'Demand for biotin'),
    @(5, 'DM_clpn140(c)', 'demand reaction for cardiolipin (tetratetradecanoyl, n-C14:0)'),
    @(6, 'DM_clpn160(c)', 'demand reaction for cardiolipin (tetrahexadecanoyl, n-C16:0)'),
    @(7, 'DM_clpn180(c)', 'demand reaction for cardiolipin (tetraoctadecanoyl, n-C18:0)'),
    @(8, 'DM_clpni16(c)', 'demand reaction for cardiolipin (14-methyl-pentadecanoyl, iso-C16)'),
    @(9, 'DM_dhptd(c)', 'Demand reaction for 4,5-dihydroxy-2,3-pentanedione'),
    @(10, 'sink_dmbzid', 'Sink for 5,6-Dimenthylbenzimidazole'),
    @(13, 'DM_4HBA', 'Sink needed to allow 4-hydroxy-benzoate to leave system'),
    @(14, 'DM_HQN', 'Sink needed to allow Hydroquinone to leave system'),
    @(15, 'DM_btn', 'Demand for biotin'),
    @(16, 'DM_clpn140(c)', 'demand reaction for cardiolipin (tetratetradecanoyl, n-C14:0)'),
    @(17, 'DM_clpn160(c)', 'demand reaction for cardiolipin (tetrahexadecanoyl, n-C16:0)'),
    @(18, 'DM_clpn180(c)', 'demand reaction for cardiolipin (tetraoctadecanoyl, n-C18:0)'),
    @(19, 'DM_clpni16(c)', 'demand reaction for cardiolipin (14-methyl-pentadecanoyl, iso-C16)'),
    @(20, 'DM_dhptd(c)', 'Demand reaction for 4,5-dihydroxy-2,3-pentanedione'),
    @(21, 'sink_dmbzid', 'Sink for 5,6-Dimenthylbenzimidazole'),
    @(24, 'DM_4HBA', 'Sink needed to allow 4-hydroxy-benzoate to leave system'),
    @(25, 'DM_5DRIB', 'Sink needed to allow 5-deoxyribose to leave system'),
    @(26, 'DM_5MTR', 'Sink needed to allow 5-Methylthio-D-ribose to leave system'),
    @(27, 'DM_GCALD', 'Sink needed to allow glycol aldehyde to leave system'),
    @(28, 'DM_HQN', 'Sink needed to allow Hydroquinone to leave system'),
    @(29, 'DM_btn', 'Demand for biotin'),
    @(30, 'DM_clpn140(c)', 'demand reaction for cardiolipin (tetratetradecanoyl, n-C14:0)'),
    @(31, 'DM_clpn160(c)', 'demand reaction for cardiolipin (tetrahexadecanoyl, n-C16:0)'),
    @(32, 'DM_clpn180(c)', 'demand reaction for cardiolipin (tetraoctadecanoyl, n-C18:0)'),
    @(33, 'DM_clpnai15(c)', 'demand reaction for cardiolipin (12-methyl-tetradecanoyl, anteiso-C15)'),
    @(34, 'DM_clpnai17(c)', 'demand reaction for cardiolipin (14-methyl-hexadecanoyl, anteiso-C17)'),
    @(35, 'DM_clpni15(c)', 'demand reaction for cardiolipin (13-methyl-tetradecanoyl, iso-C15)'),
    @(36, 'DM_clpni16(c)', 'demand reaction for cardiolipin (14-methyl-pentadecanoyl, iso-C16)'),
    @(37, 'DM_clpni17(c)', 'demand reaction for cardiolipin (15-methyl-hexadecanoyl, iso-C17)'),
    @(38, 'DM_dad_5', 'Demand for 5-deoxyadenosine'),
    @(39, 'DM_dhptd(c)', 'Demand reaction for 4,5-dihydroxy-2,3-pentanedione'),
    @(40, 'DM_teich_45_BS(c)', 'demand reaction for teichuronic acid (GlcA + GalNac, 45 repeating unit)'),
    @(41, 'sink_PGPm1[c]', 'Sink reaction for peptidoglycan polymer (n-1) subunits'),
    @(42, 'sink_dmbzid', 'Sink for 5,6-Dimenthylbenzimidazole'),
    @(43, 'sink_gthrd(c)', 'sink reaction for reduced glutathione'),
    @(46, 'DM_4HBA', 'Sink needed to allow 4-hydroxy-benzoate to leave system'),
    @(47, 'DM_5DRIB', 'Sink needed to allow 5-deoxyribose to leave system'),
    @(48, 'DM_5MTR', 'Sink needed to allow 5-Methylthio-D-ribose to leave system'),
    @(49, 'DM_GCALD', 'Sink needed to allow glycol aldehyde to leave system'),
    @(50, 'DM_HQN', 'Sink needed to allow Hydroquinone to leave system'),
    @(51, 'DM_btn', 'Demand for biotin'),
    @(52, 'DM_clpn140(c)', 'demand reaction for cardiolipin (tetratetradecanoyl, n-C14:0)'),
    @(53, 'DM_clpn160(c)', 'demand reaction for cardiolipin (tetrahexadecanoyl, n-C16:0)'),
    @(54, 'DM_clpn180(c)', 'demand reaction for cardiolipin (tetraoctadecanoyl, n-C18:0)'),
    @(55, 'DM_clpnai15(c)', 'demand reaction for cardiolipin (12-methyl-tetradecanoyl, anteiso-C15)'),
    @(56, 'DM_clpnai17(c)', 'demand reaction for cardiolipin (14-methyl-hexadecanoyl, anteiso-C17)'),
    @(57, 'DM_clpni15(c)', 'demand reaction for cardiolipin (13-methyl-tetradecanoyl, iso-C15)'),
    @(58, 'DM_clpni16(c)', 'demand reaction for cardiolipin (14-methyl-pentadecanoyl, iso-C16)'),
    @(59, 'DM_clpni17(c)', 'demand reaction for cardiolipin (15-methyl-hexadecanoyl, iso-C17)'),
    @(60, 'DM_dad_5', 'Demand for 5-deoxyadenosine'),
    @(61, 'DM_dhptd(c)', 'Demand reaction for 4,5-dihydroxy-2,3-pentanedione'),
    @(62, 'DM_teich_45_BS(c)', 'demand reaction for teichuronic acid (GlcA + GalNac, 45 repeating unit)'),
    @(63, 'sink_PGPm1[c]', 'Sink reaction for peptidoglycan polymer (n-1) subunits'),
    @(64, 'sink_gthrd(c)', 'sink reaction for reduced glutathione')
)

foreach ($h in $headers) {
    $r = $h[0]
    $ws.Cells.Item($r, 1).Value = $h[1]
    $rng = $ws.Range("A" + $r + ":D" + $r)
    $rng.HorizontalAlignment = -4131
    $rng.VerticalAlignment = -4108
    $rng.Merge()
}

foreach ($d in $data) {
    $r = $d[0]
    $ws.Cells.Item($r, 1).Value = $d[1]
    $ws.Cells.Item($r, 2).Value = $d[2]
}

# Column widths matching the source sheet
$ws.Columns.Item(1).ColumnWidth = 16.94
$ws.Columns.Item(2).ColumnWidth = 59.88

# Match the source workbook's saved selection on the new sheet
$null = $ws.Range("E26").Select()